$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.821.92"
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").Value = "2.762.33"
$ws.Range("E3").Value = "  +1.98%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "'583.99"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").Value = "'160.64"
$ws.Range("E6").Value = "  +8.10%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.618"
$ws.Range("E8").Value = "  +2.12%  "

$ws.Range("D9").Value = "2.774.00"
$ws.Range("E9").Value = "  +1.90%  "

$ws.Range("D10").Value = "'6.77"
$ws.Range("E10").Value = "  +1.64%  "

$ws.Range("D11").Value = "'0.113"
$ws.Range("E11").Value = "  +1.05%  "

$ws.Range("D12").Value = "'0.397"
$ws.Range("E12").Value = "  +2.87%  "

$ws.Range("E13").Value = "  +0.86%  "

$ws.Range("D14").Value = "3.258.38"
$ws.Range("E14").Value = "  +2.20%  "

$ws.Range("D15").Value = "'27.45"
$ws.Range("E15").Value = "  +3.85%  "

$ws.Range("D16").Value = "63.798.47"
$ws.Range("E16").Value = "  +2.05%  "

$ws.Range("D17").Value = "'0.0000158"
$ws.Range("E17").Value = "  +5.85%  "

$ws.Range("D18").Value = "2.775.02"
$ws.Range("E18").Value = "  +2.24%  "

$ws.Range("D19").Value = "'12.29"
$ws.Range("E19").Value = "  +3.54%  "

$ws.Range("D20").Value = "'4.98"
$ws.Range("E20").Value = "  +2.82%  "

$ws.Range("D21").Value = "'365.53"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("D22").Value = "'7.00"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "'0.562"
$ws.Range("E23").Value = "  +5.50%  "

$ws.Range("E24").Value = "  +0.52%  "

$ws.Range("D25").Value = "'67.28"
$ws.Range("E25").Value = "  +3.31%  "

$ws.Range("D26").Value = "'0.176"
$ws.Range("E26").Value = "  +6.65%  "

$ws.Range("D27").Value = "'8.69"
$ws.Range("E27").Value = "  +1.99%  "

$ws.Range("D28").Value = "0.0₃0958"
$ws.Range("E28").Value = "  +12.81%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("D31").Value = "'7.26"
$ws.Range("E31").Value = "  +1.74%  "

$ws.Range("E32").Value = "  +8.17%  "

$ws.Range("D33").Value = "'172.77"
$ws.Range("E33").Value = "  +1.64%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "'20.70"
$ws.Range("E35").Value = "  +1.28%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'5.04"
$ws.Range("E36").Value = "  +6.59%  "

$ws.Range("D37").Value = "'1.47"
$ws.Range("E37").Value = "  +4.61%  "

$ws.Range("D38").Value = "'1.83"
$ws.Range("E38").Value = "  +2.16%  "

$ws.Range("D39").Value = "'1.02"
$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'4.25"
$ws.Range("E40").Value = "  +1.06%  "

$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'338.61"
$ws.Range("E41").Value = "  -3.90%  "

$ws.Range("D42").Value = "'6.23"
$ws.Range("E42").Value = "  +9.46%  "

$ws.Range("D43").Value = "'39.83"
$ws.Range("E43").Value = "  +2.47%  "

$ws.Range("D44").Value = "'22.28"
$ws.Range("E44").Value = "  +4.45%  "

$ws.Range("D45").Value = "'22.66"
$ws.Range("E45").Value = "  +4.98%  "

$ws.Range("D46").Value = "'0.0607"
$ws.Range("E46").Value = "  +2.83%  "

$ws.Range("D47").Value = "'0.648"
$ws.Range("E47").Value = "  +1.59%  "

$ws.Range("D48").Value = "'0.0261"
$ws.Range("E48").Value = "  +1.53%  "

$ws.Range("D49").Value = "'138.20"
$ws.Range("E49").Value = "  +1.20%  "

$ws.Range("D50").Value = "'0.102"
$ws.Range("E50").Value = "  +1.71%  "

$ws.Range("D51").Value = "2.164.40"
$ws.Range("E51").Value = "  +1.49%  "
